$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.927.50"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "2.503.21"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.36"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.38"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.40"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.346"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "2.945.66"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").Value = "58.868.91"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.76"
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "2.511.80"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.03"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "323.28"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.93"
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.99"
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.51"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").Value = "0.0₃0760"
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.02"
$ws.Range("E29").Value = "  +1.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.45"
$ws.Range("E30").Value = "  -3.43%  "
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.34"
$ws.Range("E34").Value = "  -0.69%  "
$ws.Range("E35").Value = "  -1.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.03"
$ws.Range("E36").Value = "  -1.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.51"
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.57"
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.798"
$ws.Range("E39").Value = "  -1.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "280.46"
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.98"
$ws.Range("E43").Value = "  -5.28%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.90"
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "128.46"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0924"
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("E47").Value = "  -2.89%  "
$ws.Range("E48").Value = "  -1.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.23"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").Value = "1.749.10"
$ws.Range("E50").Value = "  -1.36%  "
$ws.Range("E51").Value = "  -0.57%  "
